$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spatial disaggregation")
$ws.Activate()

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # column F = file_path
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val.Replace("\", "/")
        $cell.Value = $newVal
    }
}

$ws.Range("F2").Select()
